# Generate-email "perfected": rewrite DoB as plain yyyy-mm-dd text and
# Email Address as Lastname-initial.Firstname@gmail.com for every student row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the DoB column to plain text so the yyyy-mm-dd strings are not
# re-interpreted as date serials when assigned below.
$ws.Range("D2:D65").NumberFormat = "@"

$ws.Range("D2").Value = "2000-06-17"
$ws.Range("E2").Value = "A.Adrian@gmail.com"
$ws.Range("D3").Value = "2001-11-14"
$ws.Range("E3").Value = "A.Shanelle@gmail.com"
$ws.Range("D4").Value = "2003-08-21"
$ws.Range("E4").Value = "A.Jude@gmail.com"
$ws.Range("D5").Value = "2003-10-26"
$ws.Range("E5").Value = "A.Kyla@gmail.com"
$ws.Range("D6").Value = "2002-03-17"
$ws.Range("E6").Value = "A.Anthony@gmail.com"
$ws.Range("D7").Value = "2002-02-24"
$ws.Range("E7").Value = "B.Moses@gmail.com"
$ws.Range("D8").Value = "2000-09-07"
$ws.Range("E8").Value = "D.Kalid@gmail.com"
$ws.Range("D9").Value = "2001-06-08"
$ws.Range("E9").Value = "E.Keith@gmail.com"
$ws.Range("D10").Value = "2001-01-28"
$ws.Range("E10").Value = "G.David@gmail.com"
$ws.Range("D11").Value = "2000-04-12"
$ws.Range("E11").Value = "G.Don@gmail.com"
$ws.Range("D12").Value = "2003-06-13"
$ws.Range("E12").Value = "H.Ashir@gmail.com"
$ws.Range("D13").Value = "2000-09-20"
$ws.Range("E13").Value = "H.Fardowsa@gmail.com"
$ws.Range("D14").Value = "2000-04-08"
$ws.Range("E14").Value = "I.Ruweida@gmail.com"
$ws.Range("D15").Value = "2001-04-27"
$ws.Range("E15").Value = "J.Myles@gmail.com"
$ws.Range("D16").Value = "2001-11-28"
$ws.Range("E16").Value = "K.Ann@gmail.com"
$ws.Range("D17").Value = "2000-02-20"
$ws.Range("E17").Value = "K.Sharon@gmail.com"
$ws.Range("D18").Value = "2001-11-24"
$ws.Range("E18").Value = "K.Neema@gmail.com"
$ws.Range("D19").Value = "2000-06-09"
$ws.Range("E19").Value = "K.Samuel@gmail.com"
$ws.Range("D20").Value = "2002-10-14"
$ws.Range("E20").Value = "K.Shannon@gmail.com"
$ws.Range("D21").Value = "2002-09-16"
$ws.Range("E21").Value = "K.Peter@gmail.com"
$ws.Range("D22").Value = "2003-10-28"
$ws.Range("E22").Value = "K.Victor@gmail.com"
$ws.Range("D23").Value = "2001-03-01"
$ws.Range("E23").Value = "K.Ian@gmail.com"
$ws.Range("D24").Value = "2000-10-02"
$ws.Range("E24").Value = "K.Eric@gmail.com"
$ws.Range("D25").Value = "2002-10-26"
$ws.Range("E25").Value = "K.Kevin@gmail.com"
$ws.Range("D26").Value = "2002-04-26"
$ws.Range("E26").Value = "K.Alex@gmail.com"
$ws.Range("D27").Value = "2000-10-23"
$ws.Range("E27").Value = "M.Kelvin@gmail.com"
$ws.Range("D28").Value = "2002-07-06"
$ws.Range("E28").Value = "M.Zivai@gmail.com"
$ws.Range("D29").Value = "2003-04-04"
$ws.Range("E29").Value = "M.Ally@gmail.com"
$ws.Range("D30").Value = "2001-02-13"
$ws.Range("E30").Value = "M.David@gmail.com"
$ws.Range("D31").Value = "2000-04-10"
$ws.Range("E31").Value = "M.Kelvin@gmail.com"
$ws.Range("D32").Value = "2000-09-26"
$ws.Range("E32").Value = "M.Victor@gmail.com"
$ws.Range("D33").Value = "2002-08-19"
$ws.Range("E33").Value = "M.Natasha@gmail.com"
$ws.Range("D34").Value = "2000-04-27"
$ws.Range("E34").Value = "M.Grace@gmail.com"
$ws.Range("D35").Value = "2001-12-28"
$ws.Range("E35").Value = "M.Mark@gmail.com"
$ws.Range("D36").Value = "2002-06-17"
$ws.Range("E36").Value = "M.Ruby@gmail.com"
$ws.Range("D37").Value = "2003-10-21"
$ws.Range("E37").Value = "M.Franklin@gmail.com"
$ws.Range("D38").Value = "2003-07-11"
$ws.Range("E38").Value = "M.Eric@gmail.com"
$ws.Range("D39").Value = "2000-08-08"
$ws.Range("E39").Value = "M.Patience@gmail.com"
$ws.Range("D40").Value = "2000-02-20"
$ws.Range("E40").Value = "M.George@gmail.com"
$ws.Range("D41").Value = "2001-07-27"
$ws.Range("E41").Value = "N.Andrew@gmail.com"
$ws.Range("D42").Value = "2002-11-28"
$ws.Range("E42").Value = "N.Monicah@gmail.com"
$ws.Range("D43").Value = "2002-07-21"
$ws.Range("E43").Value = "N.Yvonne@gmail.com"
$ws.Range("D44").Value = "2003-10-01"
$ws.Range("E44").Value = "N.Sarah@gmail.com"
$ws.Range("D45").Value = "2001-01-23"
$ws.Range("E45").Value = "N.Ian@gmail.com"
$ws.Range("D46").Value = "2002-12-22"
$ws.Range("E46").Value = "N.Alvin@gmail.com"
$ws.Range("D47").Value = "2003-08-11"
$ws.Range("E47").Value = "N.Michael@gmail.com"
$ws.Range("D48").Value = "2000-07-07"
$ws.Range("E48").Value = "N.Elizabeth@gmail.com"
$ws.Range("D49").Value = "2003-04-15"
$ws.Range("E49").Value = "N.Andrew@gmail.com"
$ws.Range("D50").Value = "2002-08-26"
$ws.Range("E50").Value = "N.Erica@gmail.com"
$ws.Range("D51").Value = "2002-10-05"
$ws.Range("E51").Value = "N.Maureen@gmail.com"
$ws.Range("D52").Value = "2000-04-02"
$ws.Range("E52").Value = "O.Steven@gmail.com"
$ws.Range("D53").Value = "2002-04-24"
$ws.Range("E53").Value = "O.Brenda@gmail.com"
$ws.Range("D54").Value = "2001-12-03"
$ws.Range("E54").Value = "O.Mwenzangu@gmail.com"
$ws.Range("D55").Value = "2002-09-27"
$ws.Range("E55").Value = "O.Caleb@gmail.com"
$ws.Range("D56").Value = "2000-01-16"
$ws.Range("E56").Value = "O.Charis@gmail.com"
$ws.Range("D57").Value = "2000-09-14"
$ws.Range("E57").Value = "O.Nicole@gmail.com"
$ws.Range("D58").Value = "2003-06-25"
$ws.Range("E58").Value = "P.Jay@gmail.com"
$ws.Range("D59").Value = "2002-06-02"
$ws.Range("E59").Value = "T.Martin@gmail.com"
$ws.Range("D60").Value = "2001-02-06"
$ws.Range("E60").Value = "T.Bramwel@gmail.com"
$ws.Range("D61").Value = "2001-08-07"
$ws.Range("E61").Value = "W.Joy@gmail.com"
$ws.Range("D62").Value = "2000-01-11"
$ws.Range("E62").Value = "W.Rosemary@gmail.com"
$ws.Range("D63").Value = "2001-04-20"
$ws.Range("E63").Value = "W.Louis@gmail.com"
$ws.Range("D64").Value = "2002-05-20"
$ws.Range("E64").Value = "W.Monika@gmail.com"
$ws.Range("D65").Value = "2002-04-19"
$ws.Range("E65").Value = "W.Trevor@gmail.com"

# Drop back to the default (unstyled) cell format now that the text values
# are safely stored, matching the original un-styled DoB column.
$ws.Range("D2:D65").Style = "Normal"

